$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (partial, in-place edits of the rich-text runs so the
# surrounding text / formatting is preserved):
#   A8: "Volume 30   Number  24"  -> "...  25"
#   C9: "Report Covering the Week  6/12/2023  Through  6/18/2023"
#       -> "...  6/19/2023  Through  6/25/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "25"
$ws.Range("C9").Characters(27, 9).Text = "6/19/2023"
$ws.Range("C9").Characters(47, 9).Text = "6/25/2023"

# ---------------------------------------------------------------------------
# Weekly crime-stat table refresh (rows 15-30).
# A few cells flip between numeric and text ("0" / "***.*") representation;
# those are written first, then their number format is copied from a
# neighbouring cell that already carries the correct style so the stored
# style index matches what a normal Excel edit would produce.
# ---------------------------------------------------------------------------

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "'***.*"
$ws.Range("F30").Value = "'0"

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -20
$ws.Range("M15").Value = 14.285714285714
$ws.Range("N15").Value = -52.941176470588

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = -24.742268041237
$ws.Range("L16").Value = 52.083333333333
$ws.Range("M16").Value = -43.846153846153
$ws.Range("N16").Value = -86.605504587156

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -15.625
$ws.Range("I17").Value = 173
$ws.Range("J17").Value = 182
$ws.Range("K17").Value = -4.945054945054
$ws.Range("L17").Value = 47.863247863247
$ws.Range("M17").Value = 67.961165048543
$ws.Range("N17").Value = 6.134969325153

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -7.692307692307
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 15.625
$ws.Range("L18").Value = 39.622641509434
$ws.Range("M18").Value = -49.31506849315
$ws.Range("N18").Value = -88.702290076335

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 149
$ws.Range("J19").Value = 172
$ws.Range("K19").Value = -13.372093023255
$ws.Range("L19").Value = 31.858407079646
$ws.Range("M19").Value = -7.453416149068
$ws.Range("N19").Value = -42.248062015503

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 11.111111111111
$ws.Range("I20").Value = 128
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = 9.401709401709
$ws.Range("L20").Value = 30.612244897959
$ws.Range("M20").Value = 9.401709401709
$ws.Range("N20").Value = -92.518994739918

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -4.040404040404
$ws.Range("I21").Value = 608
$ws.Range("J21").Value = 641
$ws.Range("K21").Value = -5.148205928237
$ws.Range("L21").Value = 37.556561085972
$ws.Range("M21").Value = -8.708708708708
$ws.Range("N21").Value = -81.86698478974

# Row 24
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -35.135135135135
$ws.Range("F24").Value = 131
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -7.746478873239
$ws.Range("I24").Value = 738
$ws.Range("J24").Value = 724
$ws.Range("K24").Value = 1.933701657458
$ws.Range("L24").Value = 53.430353430353
$ws.Range("M24").Value = 86.363636363636

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -13.207547169811
$ws.Range("I25").Value = 268
$ws.Range("J25").Value = 251
$ws.Range("K25").Value = 6.772908366533
$ws.Range("L25").Value = 34
$ws.Range("M25").Value = -10.666666666666

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 37.5
$ws.Range("L27").Value = 6.451612903225

# Row 28
$ws.Range("L28").Value = -50

# Row 29
$ws.Range("L29").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Fix up number formats / styles for the cells that switched representation
# (numeric <-> text) so the stored style matches a neighbouring cell that
# already has the right one.
# ---------------------------------------------------------------------------
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("H16").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("A26").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("A26").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("A30").Copy()
$ws.Range("F30").PasteSpecial(-4122)
